$p = $ppt.ActivePresentation

# Title slide: "Java - 3" -> "Java - 4"
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$run = $tr.Runs(5)
$run.Text = "4"
